$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 7.895089568579725
$ws.Range("E2").Value = -6.674987144438718
$ws.Range("F2").Value = -5.680488882285537
$ws.Range("G2").Value = -6.177742814835656
$ws.Range("H2").Value = 14.99878708177518
$ws.Range("I2").Value = -5.467687944985539
$ws.Range("J2").Value = 1.107023072538064
$ws.Range("D3").Value = -0.9281962351246827
$ws.Range("E3").Value = 6.72168738250815
$ws.Range("F3").Value = -0.6648559136396395
$ws.Range("G3").Value = -2.027508211865583
$ws.Range("H3").Value = 0.09793753258278383
$ws.Range("I3").Value = -0.4204894376324931
$ws.Range("J3").Value = -2.778578177568736
$ws.Range("D4").Value = 1.446086955100122
$ws.Range("E4").Value = 3.249772293330987
$ws.Range("F4").Value = -0.4595092445906481
$ws.Range("G4").Value = 0.09006123475212764
$ws.Range("H4").Value = -7.041228894227597
$ws.Range("I4").Value = 0.2885967499234481
$ws.Range("J4").Value = 2.426224687258598
$ws.Range("D5").Value = -6.599612981350536
$ws.Range("E5").Value = -0.1810729435441736
$ws.Range("F5").Value = -20.66456959606492
$ws.Range("G5").Value = 20.87879642263015
$ws.Range("H5").Value = 13.80173634841873
$ws.Range("I5").Value = 6.380599859212127
$ws.Range("J5").Value = -13.61587726669496
$ws.Range("D6").Value = -4.629890234386988
$ws.Range("E6").Value = 9.353171567705081
$ws.Range("F6").Value = 10.59986614691645
$ws.Range("G6").Value = -6.126669408505362
$ws.Range("H6").Value = 2.690645083242104
$ws.Range("I6").Value = -5.54384812702416
$ws.Range("J6").Value = -6.343301974655924
$ws.Range("D7").Value = -6.144925765434117
$ws.Range("E7").Value = 15.05711507292876
$ws.Range("F7").Value = -6.402469829866455
$ws.Range("G7").Value = 1.241502625342788
$ws.Range("H7").Value = 7.95501859872583
$ws.Range("I7").Value = -6.341082874346043
$ws.Range("J7").Value = -5.365169809854561
$ws.Range("D8").Value = -0.2161066570902882
$ws.Range("E8").Value = -1.399962587267457
$ws.Range("F8").Value = 0.4449416901078789
$ws.Range("G8").Value = 1.377899017071763
$ws.Range("H8").Value = 0.890529759708272
$ws.Range("I8").Value = 0.7764060953771168
$ws.Range("J8").Value = -1.873710373508476
$ws.Range("D9").Value = -5.161923432167462
$ws.Range("E9").Value = 3.372075722868998
$ws.Range("F9").Value = 5.06877113835789
$ws.Range("G9").Value = 2.465435716231376
$ws.Range("H9").Value = 4.203294914274939
$ws.Range("I9").Value = 1.884833668745529
$ws.Range("J9").Value = -11.83249446967147
$ws.Range("D10").Value = -3.725993835012948
$ws.Range("E10").Value = -4.183527513374042
$ws.Range("F10").Value = -4.001230094461162
$ws.Range("G10").Value = 3.060055300550586
$ws.Range("H10").Value = 9.741603874102358
$ws.Range("I10").Value = 4.07196193183082
$ws.Range("J10").Value = -4.962874015371635
